$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 09.04.2025"

$ws.Range("B6").Value = "10.04."
$ws.Range("C6").Value = "11.04."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 34598657"
$ws.Range("E6").Value = "41,80-"

$ws.Range("B7").Value = "12.04."
$ws.Range("C7").Value = "13.04."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,73-"

$ws.Range("B8").Value = "14.04."
$ws.Range("C8").Value = "15.04."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 12913608"
$ws.Range("E8").Value = "86,58-"

$ws.Range("B9").Value = "15.04."
$ws.Range("C9").Value = "16.04."
$ws.Range("D9").Value = "ZALANDO MKTPLC EU PBOHJY"
$ws.Range("E9").Value = "41,04-"

$ws.Range("D12").Value = "KONTOSTAND AM 19.04.2025"
$ws.Range("E12").Value = "194,15-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 27.04.2025"
